$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row (2-51).
# All of these were bumped by one day (45181 -> 45182), i.e. 2023-09-12 -> 2023-09-13.
$ws.Range("C2:C51").Value = 45182
